# =====================================================================
# AFter class Week 4
# Restructure DiscreateProbStarter sheet, add ContinuousProbStarter and
# ContinuousProbStarter (2) sheets, update defined names.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. DiscreateProbStarter sheet: shift the Binomial/Geometric blocks one
#    column to the right, add a new "units" column (C), and extend the
#    Poisson table down through row 25.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DiscreateProbStarter")

# Wipe out everything from column C onward so the old D/E/G/H/I layout
# doesn't leave stale cells behind once things move to E/F/G/I/J/K.
$ws.Range("C1:K25").ClearContents()

# --- Row 1 headers ---
$ws.Range("A1").Value = "Poisson"
$ws.Range("E1").Value = "Binomial"
$ws.Range("I1").Value = "Geometric"

# --- Row 2 ---
$ws.Range("E2").Value = "SingleTrialProb"
$ws.Range("F2").Value = 0.16
$ws.Range("I2").Value = "Single Trial Prob"
$ws.Range("J2").Value = 0.18

# --- Row 3 ---
$ws.Range("A3").Value = "Average Rate"
$ws.Range("B3").Value = 6.7619049999999996
$ws.Range("E3").Value = "Trials (Years etc)"
$ws.Range("F3").Value = 5

# --- Row 4 (headers for the data table) ---
$ws.Range("A4").Value = "Freq"
$ws.Range("B4").Value = "Prob"
$ws.Range("C4").Value = 21
$ws.Range("E4").Value = "Successes"
$ws.Range("F4").Value = "Prob"
$ws.Range("G4").Value = "Cum_Prob"
$ws.Range("I4").Value = "Num of Trials"
$ws.Range("J4").Value = "Prob"
$ws.Range("K4").Value = "Cum_Prob"

# --- Poisson table, rows 5-25 (A: point, B: POISSON.DIST, C: units*B) ---
$ws.Range("A5").Formula = "=0"
for ($r = 6; $r -le 25; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 5)
}
for ($r = 5; $r -le 25; $r++) {
    $ws.Cells.Item($r, 2).Formula = "=_xlfn.POISSON.DIST(A$r,averageRaber,FALSE)"
    $ws.Cells.Item($r, 3).Formula = "=units*B$r"
}

# --- Freq column D, rows 5-16 (historical raw sample data) ---
$freq = @(0,4,3,1,1,3,2,1,2,2,1,1)
for ($i = 0; $i -lt $freq.Length; $i++) {
    $ws.Cells.Item(5 + $i, 4).Value = $freq[$i]
}

# --- Binomial table, E/F/G rows 5-10 ---
for ($r = 5; $r -le 10; $r++) {
    $ws.Cells.Item($r, 5).Value = ($r - 5)
    $ws.Cells.Item($r, 6).Formula = "=_xlfn.BINOM.DIST(A$r,trials,stp,FALSE)"
    $ws.Cells.Item($r, 7).Formula = "=_xlfn.BINOM.DIST(E$r,trials,stp,TRUE)"
}

# --- Geometric table, I/J/K rows 5-10 ---
for ($r = 5; $r -le 10; $r++) {
    $ws.Cells.Item($r, 9).Value = ($r - 4)
}
$ws.Range("J5").Formula = "=gstp*((1-gstp)^(I5-1))"
$ws.Range("J6").Formula = "=gstp*((1-gstp)^(I6-1))"
$ws.Range("J7").Formula = "=gstp*((1-gstp)^(I7-1))"
$ws.Range("J8").Formula = "=gstp*((1-gstp)^(I8-1))"
$ws.Range("J9").Formula = "=gstp*((1-gstp)^(I9-1))"
$ws.Range("J10").Formula = "=gstp*((1-gstp)^(I10-1))"
$ws.Range("K5").Formula = "=J5"
$ws.Range("K6").Formula = "=K5+J6"
$ws.Range("K7").Formula = "=K6+J7"
$ws.Range("K8").Formula = "=K7+J8"
$ws.Range("K9").Formula = "=K8+J9"
$ws.Range("K10").Formula = "=K9+J10"

# Column widths follow the shifted columns (D->E, G->I)
$ws.Columns.Item(4).ColumnWidth = 8.43
$ws.Columns.Item(5).ColumnWidth = 17.3671875
$ws.Columns.Item(9).ColumnWidth = 14.1015625

$ws.Range("D17").Select()

Write-Output "DiscreateProbStarter rebuilt"

# ---------------------------------------------------------------------
# 2. Update the defined names that moved along with the columns above,
#    and add the new "units" name.
# ---------------------------------------------------------------------
$wb.Names.Item("stp").RefersTo = "=DiscreateProbStarter!`$F`$2"
$wb.Names.Item("trials").RefersTo = "=DiscreateProbStarter!`$F`$3"
$wb.Names.Item("gstp").RefersTo = "=DiscreateProbStarter!`$J`$2"
$wb.Names.Item("gt").RefersTo = "=DiscreateProbStarter!`$J`$3"
$wb.Names.Add("units", "=DiscreateProbStarter!`$C`$4")

Write-Output "defined names updated"

# ---------------------------------------------------------------------
# 3. Add the two new Continuous-probability starter sheets.
# ---------------------------------------------------------------------
$wsLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCont1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsLast)
$wsCont1.Name = "ContinuousProbStarter"

$wsCont2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsCont1)
$wsCont2.Name = "ContinuousProbStarter (2)"

function Fill-ContinuousSheet($sheet, $meanVal, $stdVal, $testVal) {
    $sheet.Range("B1").Value = "Mean"
    $sheet.Range("E1").Value = "StandardDiv"

    $sheet.Range("A3").Value = "Test Value"
    $sheet.Range("B3").Value = $meanVal
    $sheet.Range("E3").Value = "Prob"
    $sheet.Range("F3").Value = 0.9

    $sheet.Range("A4").Value = "Answer"
    $sheet.Range("B4").Value = $stdVal

    $sheet.Range("A6").Value = "Normal Find Prob"
    $sheet.Range("B6").Value = $testVal

    $sheet.Range("A8").Value = "Normal Find Cutoff"
    $sheet.Range("B8").Formula = "=_xlfn.NORM.DIST(B6,B3,B4,TRUE)"
    $sheet.Range("C8").Formula = "=1-B8"
    $sheet.Range("E8").Value = "Answer(Cutoff)"
    $sheet.Range("F8").Formula = "=_xlfn.NORM.INV(F3,B3,B4)"

    $sheet.Columns.Item(1).ColumnWidth = 15
}

Fill-ContinuousSheet $wsCont1 4695 370 5169
Fill-ContinuousSheet $wsCont2 80 20 75

# sheet-local extra rows on "ContinuousProbStarter (2)"
$wsCont2.Range("B9").Formula = "=_xlfn.NORM.DIST(90,B3,B4,TRUE)"
$wsCont2.Range("B14").Formula = "=B9-B8"

$wsCont1.Range("B6").Select()
$wsCont2.Range("B10").Select()

Write-Output "continuous sheets added"

# Sheet-scoped names on "ContinuousProbStarter (2)" plus the matching
# workbook-level names on "ContinuousProbStarter".
$wsCont2.Names.Add("mean", "='ContinuousProbStarter (2)'!`$A`$3")
$wsCont2.Names.Add("std", "='ContinuousProbStarter (2)'!`$A`$4")
$wb.Names.Add("mean", "=ContinuousProbStarter!`$A`$3")
$wb.Names.Add("std", "=ContinuousProbStarter!`$A`$4")

Write-Output "mean/std names added"

# ---------------------------------------------------------------------
# 4. Selection / active sheet bookkeeping to match the saved workbook.
# ---------------------------------------------------------------------
$wsCont2.Select()

Write-Output "done"

